$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old leftover helper cells in H1:K14 and the scratch row
# G15:K15 (single-letter codes "c","s","f","sn","gn" that aren't used
# anywhere else in the sheet any more).
$ws.Range("H1:K14").Clear()
$ws.Range("G15:K15").Clear()

# Add level 14 (row 15): same Spawners pattern as level 13.
$ws.Range("A15").Value = 14
$ws.Range("G15").Value = $ws.Range("G14").Value2

# Add level 15 (row 16): new "Win" flag in Horde(D) column and same
# Spawners pattern as level 13/14.
$ws.Range("A16").Value = 15
$ws.Range("D16").Value = "y"
$ws.Range("G16").Value = $ws.Range("G14").Value2

# Update the selection to match the author's final selection (row 13
# fully selected) and let the sheet view reflect that.
$ws.Range("A13:XFD13").Select()
